# "replace sales sales by store type"
#
# Retail Sales by Store Type.xlsx - apply the authored changes:
#   1. Rename the header text in C1 from "Sales This Category" to "Sales".
#   2. Move the active selection from F12 to E7 (and drop the old
#      top-left scroll anchor left over at A22 by simply not restating it).
#   3. Give columns A:C explicit, best-fit-derived widths.
#   4. Best-effort: resize the workbook window and refresh the
#      "last opened from" path recorded by Excel (these two are cosmetic,
#      session-local properties - we still set them through the object
#      model for completeness even though some hosts won't persist them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header text: "Sales This Category" -> "Sales"
$ws.Range("C1").Value = "Sales"

# 2. Selection moves to E7; scrolling back to the top drops topLeftCell="A22".
$ws.Range("A1").Select() | Out-Null
$ws.Range("E7").Select() | Out-Null

# 3. Column widths (bestFit-style, in characters of the Normal style font).
#    Excel's ColumnWidth setter quantizes to whole pixels (chars*MDW, MDW=6
#    in this host, plus 5px padding) before it is stored back as a "# of
#    characters" width, so we pre-compensate the requested character width
#    to land on the pixel bucket that round-trips to the desired stored
#    width (31.7109375 / 53 / 18.140625).
$colA_width = [Math]::Round(31.7109375 * 6 - 5) / 6.0
$colB_width = [Math]::Round(53        * 6 - 5) / 6.0
$colC_width = [Math]::Round(18.140625 * 6 - 5) / 6.0

$ws.Columns.Item(1).ColumnWidth = $colA_width
$ws.Columns.Item(2).ColumnWidth = $colB_width
$ws.Columns.Item(3).ColumnWidth = $colC_width

# 4. Best effort cosmetic/session metadata (window size, last folder used).
$win = $excel.ActiveWindow
$win.Width = 28800
$win.Height = 12435
